{"js": "// Fix the typo \"Giver the above tables solve the following queries.\"\n// -> \"Given the above tables solve the following queries.\"\n// (the author retyped \"Give\" + \"n\" + \" the above ...\" while correcting\n// \"Giver\" to \"Given\"; the net effect is the single-character fix below).\n\nconst body = context.document.body;\n\nconst results = body.search(\"Giver the above tables solve the following queries.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Expected text \"Giver the above tables solve the following queries.\" not found.');\n}\n\n// Replace in place so the run keeps its original formatting\n// (bold, Palatino Linotype, size 20, black themed color).\nresults.items[0].insertText(\n  \"Given the above tables solve the following queries.\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "# Fix the typo \"Giver the above tables solve the following queries.\"\n# -> \"Given the above tables solve the following queries.\"\n# (the author retyped \"Give\" + \"n\" + \" the above ...\" while correcting\n# \"Giver\" to \"Given\"; the net effect is the single-character fix below).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$find.Execute(\n    \"Giver the above tables solve the following queries.\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    $wdFindContinue,\n    $false,\n    \"Given the above tables solve the following queries.\",\n    $wdReplaceOne\n)\n"}
